$d = $word.ActiveDocument

# --- Change 1: merge the split runs in the "Articles 3-12, 12-4-6, ..." paragraph ---
# The paragraph currently reads as a single logical text spread over 5 runs; find the
# paragraph by its (already contiguous) text, then collapse the tail runs into the
# first one so only a single run/w:t remains, while keeping the original run's rPr.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "Articles 3-12, 12-4-6, 14, 15 et 27 C. civ.") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $pStart = $target.Range.Start
    $pEnd = $target.Range.End
    $headLen = 11  # length of "Articles 3-"
    $rest = $d.Range($pStart + $headLen, $pEnd - 1)
    $rest.Delete()
    $insertPoint = $d.Range($pStart + $headLen, $pStart + $headLen)
    $insertPoint.InsertAfter("12, 12-4-6, 14, 15 et 27 C. civ.")
}

# --- Change 2: append new paragraphs after the final paragraph in the document ---
for ($n = 0; $n -lt 6; $n++) {
    $last = $d.Paragraphs.Item($d.Paragraphs.Count)
    $last.Range.InsertParagraphAfter()
}

$base = $d.Paragraphs.Count - 6
$d.Paragraphs.Item($base + 2).Range.InsertAfter("Art. L. 112-1 C. pr. Int. Article L. 331-24 du CPI.")
$d.Paragraphs.Item($base + 4).Range.InsertAfter("Articles 131-4 et 225-7-1 c. pén.")
$d.Paragraphs.Item($base + 6).Range.InsertAfter("Art. 694-4-1 et ")
